# Auto-generated edit script: updates crypto price/volume table cells
# per the target diff (Sat Apr 27 07:47:00 UTC 2024 GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.097.88"
$ws.Range("E2").Value = "  -1.98%  "
# Row 3
$ws.Range("D3").Value = "3.127.28"
$ws.Range("E3").Value = "  -0.35%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.63"
$ws.Range("E5").Value = "  -2.45%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.47"
$ws.Range("E6").Value = "  -5.25%  "
# Row 7
$ws.Range("E7").Value = "  +0.11%  "
# Row 8
$ws.Range("D8").Value = "3.121.13"
$ws.Range("E8").Value = "  -0.37%  "
# Row 9
$ws.Range("E9").Value = "  -2.76%  "
# Row 10
$ws.Range("E10").Value = "  -3.59%  "
# Row 11
$ws.Range("E11").Value = "  -2.38%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  -3.36%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("E13").Value = "  -3.36%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.15"
$ws.Range("E14").Value = "  -3.75%  "
# Row 15
$ws.Range("D15").Value = "3.642.30"
$ws.Range("E15").Value = "  -0.28%  "
# Row 16
$ws.Range("E16").Value = "  +3.24%  "
# Row 17
$ws.Range("D17").Value = "63.123.55"
$ws.Range("E17").Value = "  -1.87%  "
# Row 18
$ws.Range("D18").Value = "3.127.31"
$ws.Range("E18").Value = "  -1.34%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.67"
$ws.Range("E19").Value = "  -2.80%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "472.16"
$ws.Range("E20").Value = "  -1.18%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.16"
$ws.Range("E21").Value = "  -4.95%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  -3.20%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.68"
$ws.Range("E23").Value = "  -1.51%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "86.58"
$ws.Range("E24").Value = "  +1.12%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.95"
$ws.Range("E25").Value = "  -3.98%  "
# Row 26
$ws.Range("E26").Value = "  -0.04%  "
# Row 27
$ws.Range("E27").Value = "  -1.91%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.12"
$ws.Range("E28").Value = "  -2.41%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.94"
$ws.Range("E29").Value = "  -6.52%  "
# Row 30
$ws.Range("E30").Value = "  -0.36%  "
# Row 31
$ws.Range("E31").Value = "  +0.15%  "
# Row 32
$ws.Range("E32").Value = "  -1.53%  "
# Row 33
$ws.Range("E33").Value = "  -6.45%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("E34").Value = "  -4.61%  "
# Row 35
$ws.Range("E35").Value = "  -2.63%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").Value = "  -2.93%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.03"
$ws.Range("E37").Value = "  -0.97%  "
# Row 38
$ws.Range("D38").Value = "0.0₃0710"
$ws.Range("E38").Value = "  -5.19%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "424.20"
$ws.Range("E39").Value = "  -5.12%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0386"
$ws.Range("E40").Value = "  -1.98%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.22"
$ws.Range("E41").Value = "  -1.04%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("E42").Value = "  -10.26%  "
# Row 43
$ws.Range("D43").Value = "2.896.24"
$ws.Range("E43").Value = "  +0.09%  "
# Row 44
$ws.Range("E44").Value = "  -4.67%  "
# Row 45
$ws.Range("E45").Value = "  -0.35%  "
# Row 47
$ws.Range("E47").Value = "  -5.26%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.58"
$ws.Range("E48").Value = "  -2.93%  "
# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.112"
$ws.Range("E49").Value = "  -1.11%  "
# Row 50
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.27"
$ws.Range("E50").Value = "  -6.02%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.50"
$ws.Range("E51").Value = "  +0.15%  "
